$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "<?xml version=`"1.0`" encoding=`"UTF-8`"?>
<rpc-reply message-id=`"urn:uuid:e97bf60f-3e4b-4d13-9bb8-52a57e15824e`"
 xmlns:ncx=`"http://netconfcentral.org/ns/yuma-ncx`"
 ncx:last-modified=`"2020-10-07T13:51:28Z`" ncx:etag=`"814`"
 xmlns=`"urn:ietf:params:xml:ns:netconf:base:1.0`">
 <data>
  <components xmlns=`"http://openconfig.net/yang/platform`">
   <component>
    <name>Waveserver-Ai</name>
    <state>
     <description>Waveserver Ai Chassis 3-slot, 1RU</description>
    </state>
   </component>
  </components>
 </data>
</rpc-reply>"

$ws.Range("J3").Value = "<?xml version=`"1.0`" encoding=`"UTF-8`"?>
<rpc-reply message-id=`"urn:uuid:46ab6473-6fbe-44da-be2e-1a95719f1517`"
 xmlns:ncx=`"http://netconfcentral.org/ns/yuma-ncx`"
 ncx:last-modified=`"2020-10-07T13:51:28Z`" ncx:etag=`"814`"
 xmlns=`"urn:ietf:params:xml:ns:netconf:base:1.0`">
 <data>
  <components xmlns=`"http://openconfig.net/yang/platform`">
   <component>
    <name>Waveserver-Ai</name>
    <state>
     <hardware-version>001</hardware-version>
    </state>
   </component>
  </components>
 </data>
</rpc-reply>"

$ws.Range("J4").Value = "<?xml version=`"1.0`" encoding=`"UTF-8`"?>
<rpc-reply message-id=`"urn:uuid:65ca91d0-06b5-4aa4-b2e8-8d8184255e65`"
 xmlns:ncx=`"http://netconfcentral.org/ns/yuma-ncx`"
 ncx:last-modified=`"2020-10-07T13:51:28Z`" ncx:etag=`"814`"
 xmlns=`"urn:ietf:params:xml:ns:netconf:base:1.0`">
 <data>
  <components xmlns=`"http://openconfig.net/yang/platform`">
   <component>
    <name>Waveserver-Ai</name>
    <state>
     <id>Waveserver Ai Chassis</id>
    </state>
   </component>
  </components>
 </data>
</rpc-reply>"

$ws.Range("J5").Value = "<?xml version=`"1.0`" encoding=`"UTF-8`"?>
<rpc-reply message-id=`"urn:uuid:b11e67a3-bddb-42b5-bcc1-1e9d4ebb6ed2`"
 xmlns:ncx=`"http://netconfcentral.org/ns/yuma-ncx`"
 ncx:last-modified=`"2020-10-07T13:51:28Z`" ncx:etag=`"814`"
 xmlns=`"urn:ietf:params:xml:ns:netconf:base:1.0`">
 <data>
  <components xmlns=`"http://openconfig.net/yang/platform`">
   <component>
    <name>Waveserver-Ai</name>
    <state>
    </state>
   </component>
  </components>
 </data>
</rpc-reply>"
